$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 21052
$ws.Range("I93").Value = 21052
$ws.Range("K93").Value = 21052
$ws.Range("M93").Value = -18556
$ws.Range("H115").Value = 1056.9166
$ws.Range("I115").Value = 645.75
$ws.Range("K115").Value = 1937.25
$ws.Range("M115").Value = -370.25
$ws.Range("H116").Value = 7744.8696
$ws.Range("I116").Value = 2666
$ws.Range("J116").Value = 11651.692
$ws.Range("K116").Value = 2666
$ws.Range("L116").Value = 11651.692
$ws.Range("M116").Value = 776
$ws.Range("N116").Value = -18535.692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 24447
$ws.Range("J51").Value = 24447
$ws.Range("L51").Value = 24447
$ws.Range("N51").Value = -25959
$ws.Range("H103").Value = 38000
$ws.Range("J103").Value = 38000
$ws.Range("L103").Value = 38000
$ws.Range("N103").Value = -40344

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1268.4412
$ws.Range("I134").Value = 1047.5927
$ws.Range("J134").Value = 2120.2856
$ws.Range("K134").Value = 3142.7781
$ws.Range("L134").Value = 6360.8568
$ws.Range("M134").Value = -607.7780999999995
$ws.Range("N134").Value = -11430.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3691.2341
$ws.Range("I31").Value = 1501.56
$ws.Range("J31").Value = 4484.594
$ws.Range("K31").Value = 1501.56
$ws.Range("L31").Value = 4484.594
$ws.Range("M31").Value = -1206.56
$ws.Range("N31").Value = -5074.594
$ws.Range("H34").Value = 3691.2341
$ws.Range("I34").Value = 1501.56
$ws.Range("J34").Value = 4484.594
$ws.Range("K34").Value = 1501.56
$ws.Range("L34").Value = 4484.594
$ws.Range("M34").Value = -1299.56
$ws.Range("N34").Value = -4888.594
$ws.Range("H43").Value = 16811.111
$ws.Range("J43").Value = 16811.111
$ws.Range("L43").Value = 16811.111
$ws.Range("N43").Value = -17179.111
$ws.Range("H95").Value = 14008
$ws.Range("J95").Value = 14008
$ws.Range("L95").Value = 14008
$ws.Range("N95").Value = -19500
$ws.Range("H101").Value = 16811.111
$ws.Range("J101").Value = 16811.111
$ws.Range("L101").Value = 16811.111
$ws.Range("N101").Value = -23301.111

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 684.1818
$ws.Range("I5").Value = 597.4375
$ws.Range("J5").Value = 915.5
$ws.Range("K5").Value = 1792.3125
$ws.Range("L5").Value = 2746.5
$ws.Range("M5").Value = -1680.3125
$ws.Range("N5").Value = -2970.5
$ws.Range("H62").Value = 2913.2856
$ws.Range("J62").Value = 2913.2856
$ws.Range("L62").Value = 8739.856800000001
$ws.Range("N62").Value = -10111.8568
$ws.Range("H65").Value = 2913.2856
$ws.Range("J65").Value = 2913.2856
$ws.Range("L65").Value = 26219.5704
$ws.Range("N65").Value = -33083.5704
$ws.Range("H75").Value = 8500
$ws.Range("J75").Value = 8500
$ws.Range("L75").Value = 25500
$ws.Range("N75").Value = -27496
$ws.Range("H78").Value = 8500
$ws.Range("J78").Value = 8500
$ws.Range("L78").Value = 76500
$ws.Range("N78").Value = -86484
$ws.Range("H92").Value = 851.5
$ws.Range("J92").Value = 851.5
$ws.Range("L92").Value = 2554.5
$ws.Range("N92").Value = -5050.5
$ws.Range("H97").Value = 443.5
$ws.Range("I97").Value = 324.33334
$ws.Range("J97").Value = 483.22223
$ws.Range("K97").Value = 973.0000200000001
$ws.Range("L97").Value = 1449.66669
$ws.Range("M97").Value = -477.0000200000001
$ws.Range("N97").Value = -2441.66669
$ws.Range("H98").Value = 392.30768
$ws.Range("I98").Value = 300
$ws.Range("J98").Value = 900
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 2700
$ws.Range("M98").Value = 598
$ws.Range("N98").Value = -5696
$ws.Range("H107").Value = 673
$ws.Range("J107").Value = 1296
$ws.Range("L107").Value = 3888
$ws.Range("N107").Value = -7728
$ws.Range("H113").Value = 788103.9
$ws.Range("I113").Value = 464.3
$ws.Range("J113").Value = 1181923.6
$ws.Range("K113").Value = 1392.9
$ws.Range("L113").Value = 3545770.8
$ws.Range("M113").Value = 777.0999999999999
$ws.Range("N113").Value = -3550110.8
$ws.Range("H122").Value = 595.2778
$ws.Range("I122").Value = 238
$ws.Range("J122").Value = 952.55554
$ws.Range("K122").Value = 2142
$ws.Range("L122").Value = 8572.99986
$ws.Range("M122").Value = 308
$ws.Range("N122").Value = -13472.99986
$ws.Range("H132").Value = 3258742.8
$ws.Range("I132").Value = 4762410
$ws.Range("J132").Value = 101041.5
$ws.Range("K132").Value = 42861690
$ws.Range("L132").Value = 909373.5
$ws.Range("M132").Value = -42859160
$ws.Range("N132").Value = -914433.5
$ws.Range("H135").Value = 684.1818
$ws.Range("I135").Value = 597.4375
$ws.Range("J135").Value = 915.5
$ws.Range("K135").Value = 5376.9375
$ws.Range("L135").Value = 8239.5
$ws.Range("M135").Value = -2841.9375
$ws.Range("N135").Value = -13309.5
$ws.Range("H138").Value = 2130.818
$ws.Range("I138").Value = 1942.5
$ws.Range("J138").Value = 2633
$ws.Range("K138").Value = 5827.5
$ws.Range("L138").Value = 7899
$ws.Range("M138").Value = -687.5
$ws.Range("N138").Value = -18179
$ws.Range("H141").Value = 74167.36
$ws.Range("I141").Value = 93040.27
$ws.Range("J141").Value = 4966.6665
$ws.Range("K141").Value = 279120.81
$ws.Range("L141").Value = 14899.9995
$ws.Range("M141").Value = -273940.81
$ws.Range("N141").Value = -25259.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10013.818
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 10925.2
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 10925.2
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -11515.2
$ws.Range("H27").Value = 10013.818
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 10925.2
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 10925.2
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -11139.2
$ws.Range("H46").Value = 2943.8333
$ws.Range("J46").Value = 2167.3333
$ws.Range("L46").Value = 2167.3333
$ws.Range("N46").Value = -2543.3333
$ws.Range("H95").Value = 24158.8
$ws.Range("J95").Value = 24158.8
$ws.Range("L95").Value = 24158.8
$ws.Range("N95").Value = -29650.8
$ws.Range("H103").Value = 19950.75
$ws.Range("J103").Value = 19950.75
$ws.Range("L103").Value = 19950.75
$ws.Range("N103").Value = -22294.75
